# This workbook's rows 2-15 got reordered (a weekly refresh of the same
# logical rows in a different row order) while the header row (row 1) and
# the columns that are identical across every row (A, B, C, E, F, G, H, I,
# J, K, R, T) stay untouched. Apply the new row order by writing explicit
# values into the columns that actually vary: D (Fecha), L (Calidad),
# M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), Q (Unidad de comercializacion) and
# S (Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per destination row (2..15), taken from the rest of the sheet
# per the target layout.
$rows = @{
    2  = @{ D = 44819; L = "Primera"; M = 100; N = 25000; O = 26000; P = 25500; Q = "`$/caja 18 kilos granel"; S = 1417 }
    3  = @{ D = 44316; L = "Primera"; M = 50;  N = 20000; O = 20000; P = 20000; Q = "`$/caja 18 kilos";        S = 1111 }
    4  = @{ D = 45084; L = "Primera"; M = 100; N = 20000; O = 21000; P = 20500; Q = "`$/caja 18 kilos granel"; S = 1139 }
    5  = @{ D = 44687; L = "Primera"; M = 100; N = 18000; O = 19000; P = 18500; Q = "`$/caja 18 kilos";        S = 1028 }
    6  = @{ D = 45044; L = "Primera"; M = 100; N = 17000; O = 18000; P = 17500; Q = "`$/caja 18 kilos";        S = 972  }
    7  = @{ D = 45030; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; S = 861  }
    8  = @{ D = 44516; L = "Primera"; M = 100; N = 33000; O = 34000; P = 33500; Q = "`$/caja 18 kilos";        S = 1861 }
    9  = @{ D = 45002; L = "Primera"; M = 100; N = 12000; O = 13000; P = 12500; Q = "`$/caja 18 kilos";        S = 694  }
    10 = @{ D = 45014; L = "Primera"; M = 50;  N = 13000; O = 14000; P = 13600; Q = "`$/caja 18 kilos";        S = 756  }
    11 = @{ D = 45014; L = "Segunda"; M = 20;  N = 10000; O = 10000; P = 10000; Q = "`$/caja 18 kilos";        S = 556  }
    12 = @{ D = 44280; L = "Primera"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos";        S = 806  }
    13 = @{ D = 44280; L = "Segunda"; M = 50;  N = 12000; O = 12000; P = 12000; Q = "`$/caja 18 kilos";        S = 667  }
    14 = @{ D = 44699; L = "Primera"; M = 100; N = 20000; O = 22000; P = 21000; Q = "`$/caja 18 kilos";        S = 1167 }
    15 = @{ D = 44699; L = "Segunda"; M = 50;  N = 18000; O = 18000; P = 18000; Q = "`$/caja 18 kilos";        S = 1000 }
}

foreach ($r in $rows.Keys) {
    $row = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $row.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $row.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $row.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $row.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $row.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $row.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $row.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $row.S   # S: Precio $/Kg
}
